$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 0.945207547255415
$ws.Range("D4").Value = 4.22977938105215
$ws.Range("D5").Value = 14.479566725029
$ws.Range("D6").Value = 36.2914380824517
$ws.Range("D7").Value = 71.7159972479564
